$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.36219633333334
$ws.Range("H2").Value = 37.086589
$ws.Range("I2").Value = 0.8692805094072583
$ws.Range("J2").Value = 0.8692805094072584
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 33.380049
$ws.Range("N2").Value = 100.140147
$ws.Range("O2").Value = 0.3891462059670435
$ws.Range("P2").Value = 0.3891462059670435
$ws.Range("Q2").Value = 412.650719354287
$ws.Range("R2").Value = 3713.856474188583
$ws.Range("S2").Value = 0.3382772121569335
$ws.Range("T2").Value = 0.3382772121569335

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.36219633333334
$ws.Range("H3").Value = 37.086589
$ws.Range("I3").Value = 0.8692805094072583
$ws.Range("J3").Value = 0.8692805094072584
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 33.85786133333334
$ws.Range("N3").Value = 101.573584
$ws.Range("O3").Value = 0.3947165649764305
$ws.Range("P3").Value = 0.3947165649764305
$ws.Range("Q3").Value = 418.5575292294419
$ws.Range("R3").Value = 3767.017763064977
$ws.Range("S3").Value = 0.3431194166741947
$ws.Range("T3").Value = 0.3431194166741947

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.36219633333334
$ws.Range("H4").Value = 37.086589
$ws.Range("I4").Value = 0.8692805094072583
$ws.Range("J4").Value = 0.8692805094072584
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.53974466666667
$ws.Range("N4").Value = 55.61923400000001
$ws.Range("O4").Value = 0.2161372290565261
$ws.Range("P4").Value = 0.2161372290565261
$ws.Range("Q4").Value = 229.1919635392029
$ws.Range("R4").Value = 2062.727671852826
$ws.Range("S4").Value = 0.1878838805761303
$ws.Range("T4").Value = 0.1878838805761303

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.434409
$ws.Range("H5").Value = 4.303227
$ws.Range("I5").Value = 0.1008642600875229
$ws.Range("J5").Value = 0.1008642600875229
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.380049
$ws.Range("N5").Value = 100.140147
$ws.Range("O5").Value = 0.3891462059670435
$ws.Range("P5").Value = 0.3891462059670435
$ws.Range("Q5").Value = 47.88064270604099
$ws.Range("R5").Value = 430.9257843543689
$ws.Range("S5").Value = 0.03925094413073264
$ws.Range("T5").Value = 0.03925094413073265

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.434409
$ws.Range("H6").Value = 4.303227
$ws.Range("I6").Value = 0.1008642600875229
$ws.Range("J6").Value = 0.1008642600875229
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 33.85786133333334
$ws.Range("N6").Value = 101.573584
$ws.Range("O6").Value = 0.3947165649764305
$ws.Range("P6").Value = 0.3947165649764305
$ws.Range("Q6").Value = 48.56602101728534
$ws.Range("R6").Value = 437.094189155568
$ws.Range("S6").Value = 0.03981279427063633
$ws.Range("T6").Value = 0.03981279427063633

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.434409
$ws.Range("H7").Value = 4.303227
$ws.Range("I7").Value = 0.1008642600875229
$ws.Range("J7").Value = 0.1008642600875229
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.53974466666667
$ws.Range("N7").Value = 55.61923400000001
$ws.Range("O7").Value = 0.2161372290565261
$ws.Range("P7").Value = 0.2161372290565261
$ws.Range("Q7").Value = 26.59357660756866
$ws.Range("R7").Value = 239.342189468118
$ws.Range("S7").Value = 0.02180052168615396
$ws.Range("T7").Value = 0.02180052168615397

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4245766666666667
$ws.Range("H8").Value = 1.27373
$ws.Range("I8").Value = 0.02985523050521867
$ws.Range("J8").Value = 0.02985523050521867
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 33.380049
$ws.Range("N8").Value = 100.140147
$ws.Range("O8").Value = 0.3891462059670435
$ws.Range("P8").Value = 0.3891462059670435
$ws.Range("Q8").Value = 14.17238993759
$ws.Range("R8").Value = 127.55150943831
$ws.Range("S8").Value = 0.01161804967937738
$ws.Range("T8").Value = 0.01161804967937738

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4245766666666667
$ws.Range("H9").Value = 1.27373
$ws.Range("I9").Value = 0.02985523050521867
$ws.Range("J9").Value = 0.02985523050521867
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.85786133333334
$ws.Range("N9").Value = 101.573584
$ws.Range("O9").Value = 0.3947165649764305
$ws.Range("P9").Value = 0.3947165649764305
$ws.Range("Q9").Value = 14.37525790536889
$ws.Range("R9").Value = 129.37732114832
$ws.Range("S9").Value = 0.01178435403159945
$ws.Range("T9").Value = 0.01178435403159945

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4245766666666667
$ws.Range("H10").Value = 1.27373
$ws.Range("I10").Value = 0.02985523050521867
$ws.Range("J10").Value = 0.02985523050521867
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 18.53974466666667
$ws.Range("N10").Value = 55.61923400000001
$ws.Range("O10").Value = 0.2161372290565261
$ws.Range("P10").Value = 0.2161372290565261
$ws.Range("Q10").Value = 7.871542991424445
$ws.Range("R10").Value = 70.84388692282
$ws.Range("S10").Value = 0.006452826794241831
$ws.Range("T10").Value = 0.006452826794241832

Write-Output "Updated rows 2-10 with new NATMI edge statistics."
